$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K = strikeouts), rows 2-28, replacing old Strike# values
$gValues = @{
    2  = 3
    3  = 4
    4  = 0
    5  = 2
    6  = 11
    7  = 4
    8  = 3
    9  = 4
    10 = 2
    11 = 4
    12 = 5
    13 = 6
    14 = 2
    15 = 4
    16 = 0
    17 = 8
    18 = 4
    19 = 6
    20 = 3
    21 = 7
    22 = 5
    23 = 3
    24 = 7
    25 = 5
    26 = 3
    27 = 1
    28 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
